# S7 Symbols.xlsx: added Unixtime and UDT_NET_CONFIG symbols
#
# The symbol table on "Tabelle1" lists Name / Address / Address / Comment
# columns. Two new rows are inserted:
#   - "Unixtime"       / "FC     112"  -> inserted as the new row 22
#   - "UDT_NET_CONFIG" / "UDT     101" -> appended as the new row 24
# (the previous last row, UDT_STRUCTANYPTR / UDT     100, shifts from
# row 22 down to row 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Make room for the new "Unixtime" row by inserting a blank row at 22 -
# this shifts the old row 22 (UDT_STRUCTANYPTR) down to row 23 and carries
# its formatting (style "locked=0") along for the ride.
$ws.Rows("22:22").Insert()

# New last row (24): UDT_NET_CONFIG / UDT     101 / UDT     101
$ws.Range("A24").Value = "UDT_NET_CONFIG"
$ws.Range("B24").Value = "UDT     101"
$ws.Range("C24").Value = "UDT     101"
$ws.Range("B24:C24").Locked = $false

# New row 22: Unixtime / FC     112 / FC     112
$ws.Range("A22").Value = "Unixtime"
$ws.Range("B22").Value = "FC     112"
$ws.Range("C22").Value = "FC     112"

# Re-apply the sheet protection present in the edited workbook.
$ws.Protect("840F", $true, $true, $true)
